$d = $word.ActiveDocument

# Helper: pin a run boundary at the edges of a Range by toggling a
# formatting property on and back off. This forces the text inside the
# range to live in its own run, separate from neighboring text, even
# though the visible formatting ends up identical again.
function Pin-RunBoundary($range) {
    $range.Font.Bold = 1
    $range.Font.Bold = 0
}

# ---------------------------------------------------------------------------
# Change 1: paragraph "Министерство образования Московской области"
# becomes two runs: "Министерство образования " + "ХХХ"
# (i.e. "Московской области" is replaced with "ХХХ")
# ---------------------------------------------------------------------------
$rng1 = $d.Content
$found1 = $rng1.Find.Execute("Министерство образования Московской области", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found1) {
    throw "Could not find target text for change 1"
}

$prefix1 = "Министерство образования "
$oldTail1 = "Московской области"
$tailStart1 = $rng1.Start + $prefix1.Length
$tailEnd1 = $tailStart1 + $oldTail1.Length
$target1 = $d.Range($tailStart1, $tailEnd1)

$target1.Delete()
$target1.Collapse(0)
$target1.InsertAfter("ХХХ")

$newRng1 = $d.Range($target1.Start, $target1.Start + 3)
Pin-RunBoundary $newRng1

# ---------------------------------------------------------------------------
# Change 2: paragraph starting with "Московской области «ХХХХ»"
# becomes: "ХХХ" + " «" + "ХХХХ" + "»"
# (i.e. "Московской области" is replaced with "ХХХ"; the rest of the
# paragraph is untouched)
# ---------------------------------------------------------------------------
$rng2 = $d.Content
$found2 = $rng2.Find.Execute("Московской области «", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found2) {
    throw "Could not find target text for change 2"
}

$oldHead2 = "Московской области"
$headStart2 = $rng2.Start
$headEnd2 = $headStart2 + $oldHead2.Length
$target2 = $d.Range($headStart2, $headEnd2)

$target2.Delete()
$target2.Collapse(0)
$target2.InsertAfter("ХХХ")

$newRng2 = $d.Range($target2.Start, $target2.Start + 3)
Pin-RunBoundary $newRng2

# The leftover " «" text remains right after the newly inserted "ХХХ";
# pin its trailing boundary too so it doesn't get merged with the
# following, already-existing "ХХХХ" run.
$suffix2 = " «"
$suffixRng2 = $d.Range($newRng2.End, $newRng2.End + $suffix2.Length)
Pin-RunBoundary $suffixRng2

# The existing "ХХХХ" run (followed by the existing "»" run) must also
# keep its own boundary -- otherwise it gets swept into the run we just
# touched above. Pin it too, without altering its text.
$xxxx2 = "ХХХХ"
$xxxxRng2 = $d.Range($suffixRng2.End, $suffixRng2.End + $xxxx2.Length)
Pin-RunBoundary $xxxxRng2

Write-Host "Change 1 applied: $found1"
Write-Host "Change 2 applied: $found2"
